$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.223.51'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.29'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.14'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4713'
$ws.Range("E7").Value = '  +1.00%  '
$ws.Range("E8").Value = '  +2.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06567'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.86'
$ws.Range("E10").Value = '  +3.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07941'
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.76'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.857.80'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.135'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6808'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '267.15'
$ws.Range("E16").Value = '  -4.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.225.45'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.69'
$ws.Range("E18").Value = '  +8.28%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007542'
$ws.Range("E19").Value = '  +3.57%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.101.20'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.269'
$ws.Range("E23").Value = '  -4.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.170'
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.58'
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.180'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.951'
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.399'
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09940'
$ws.Range("E30").Value = '  +2.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.335'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.470'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.013'
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04705'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7013'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.623'
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.325'
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.00'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.943'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8415'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4161'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.13'
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.151'
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '944.19'
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.230'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.14'
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05659'
$ws.Range("E51").Value = '  +0.64%  '
